$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2851.4
$ws.Range("I62").Value = 3009.25
$ws.Range("K62").Value = 3009.25
$ws.Range("M62").Value = -2385.25
$ws.Range("H65").Value = 2851.4
$ws.Range("I65").Value = 3009.25
$ws.Range("K65").Value = 15046.25
$ws.Range("M65").Value = -11926.25
$ws.Range("H80").Value = 4248.069
$ws.Range("I80").Value = 402.66666
$ws.Range("J80").Value = 6962.4707
$ws.Range("K80").Value = 1207.99998
$ws.Range("L80").Value = 20887.4121
$ws.Range("M80").Value = -209.9999800000001
$ws.Range("N80").Value = -22883.4121
$ws.Range("H83").Value = 4248.069
$ws.Range("I83").Value = 402.66666
$ws.Range("J83").Value = 6962.4707
$ws.Range("K83").Value = 3623.99994
$ws.Range("L83").Value = 62662.2363
$ws.Range("M83").Value = 1368.00006
$ws.Range("N83").Value = -72646.23629999999
$ws.Range("H137").Value = 1801.0358
$ws.Range("I137").Value = 1354.9269
$ws.Range("J137").Value = 3020.4
$ws.Range("K137").Value = 4064.7807
$ws.Range("L137").Value = 9061.200000000001
$ws.Range("M137").Value = -1514.7807
$ws.Range("N137").Value = -14161.2
$ws.Range("H138").Value = 3173.4773
$ws.Range("I138").Value = 1443.7097
$ws.Range("J138").Value = 7298.3076
$ws.Range("K138").Value = 4331.1291
$ws.Range("L138").Value = 21894.9228
$ws.Range("M138").Value = 808.8708999999999
$ws.Range("N138").Value = -32174.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1158.5
$ws.Range("I2").Value = 1323.25
$ws.Range("K2").Value = 1323.25
$ws.Range("M2").Value = -1210.25
$ws.Range("H116").Value = 1158.5
$ws.Range("I116").Value = 1323.25
$ws.Range("K116").Value = 1323.25
$ws.Range("M116").Value = 970.75
$ws.Range("H132").Value = 4150.362
$ws.Range("I132").Value = 1358.2759
$ws.Range("J132").Value = 8648.723
$ws.Range("K132").Value = 4074.8277
$ws.Range("L132").Value = 25946.169
$ws.Range("M132").Value = -1544.8277
$ws.Range("N132").Value = -31006.169

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1158.5
$ws.Range("I3").Value = 1323.25
$ws.Range("K3").Value = 1323.25
$ws.Range("M3").Value = -1209.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2146.9583
$ws.Range("I31").Value = 1572.4286
$ws.Range("J31").Value = 6168.6665
$ws.Range("K31").Value = 1572.4286
$ws.Range("L31").Value = 6168.6665
$ws.Range("M31").Value = -1277.4286
$ws.Range("N31").Value = -6758.6665
$ws.Range("H34").Value = 2146.9583
$ws.Range("I34").Value = 1572.4286
$ws.Range("J34").Value = 6168.6665
$ws.Range("K34").Value = 1572.4286
$ws.Range("L34").Value = 6168.6665
$ws.Range("M34").Value = -1370.4286
$ws.Range("N34").Value = -6572.6665
$ws.Range("H58").Value = 1282110.9
$ws.Range("I58").Value = 1624458.5
$ws.Range("K58").Value = 1624458.5
$ws.Range("M58").Value = -1624255.5
$ws.Range("H136").Value = 1282110.9
$ws.Range("I136").Value = 1624458.5
$ws.Range("K136").Value = 4873375.5
$ws.Range("M136").Value = -4870825.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2116
$ws.Range("I34").Value = 239.75
$ws.Range("J34").Value = 3054.125
$ws.Range("K34").Value = 719.25
$ws.Range("L34").Value = 9162.375
$ws.Range("M34").Value = -635.25
$ws.Range("N34").Value = -9330.375
$ws.Range("H39").Value = 9890.4375
$ws.Range("J39").Value = 9890.4375
$ws.Range("L39").Value = 29671.3125
$ws.Range("N39").Value = -30259.3125
$ws.Range("H55").Value = 3761.3076
$ws.Range("I55").Value = 1230.6666
$ws.Range("J55").Value = 4520.5
$ws.Range("K55").Value = 3691.9998
$ws.Range("L55").Value = 13561.5
$ws.Range("M55").Value = -3514.9998
$ws.Range("N55").Value = -13915.5
$ws.Range("H131").Value = 1700.5
$ws.Range("J131").Value = 1436.0416
$ws.Range("L131").Value = 4308.1248
$ws.Range("N131").Value = -14388.1248
$ws.Range("H132").Value = 1812.5834
$ws.Range("I132").Value = 1995.3
$ws.Range("J132").Value = 1682.0714
$ws.Range("K132").Value = 17957.7
$ws.Range("L132").Value = 15138.6426
$ws.Range("M132").Value = -15427.7
$ws.Range("N132").Value = -20198.6426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3857.1428
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -2502
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 3857.1428
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 17500
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -12508
$ws.Range("N83").Value = -29984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6426.5586
$ws.Range("I7").Value = 4876.143
$ws.Range("J7").Value = 8931.076999999999
$ws.Range("K7").Value = 4876.143
$ws.Range("L7").Value = 8931.076999999999
$ws.Range("M7").Value = -4764.143
$ws.Range("N7").Value = -9155.076999999999
$ws.Range("H68").Value = 1000
$ws.Range("I68").Value = 1000
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 1000
$ws.Range("L68").Value = 1000
$ws.Range("M68").Value = -251
$ws.Range("N68").Value = -2498
$ws.Range("H71").Value = 1000
$ws.Range("I71").Value = 1000
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 5000
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = -1256
$ws.Range("N71").Value = -12488
$ws.Range("H126").Value = 6426.5586
$ws.Range("I126").Value = 4876.143
$ws.Range("J126").Value = 8931.076999999999
$ws.Range("K126").Value = 14628.429
$ws.Range("L126").Value = 26793.231
$ws.Range("M126").Value = -12158.429
$ws.Range("N126").Value = -31733.231
$ws.Range("H136").Value = 3181.7466
$ws.Range("I136").Value = 1704.5862
$ws.Range("J136").Value = 8221.471
$ws.Range("K136").Value = 5113.7586
$ws.Range("L136").Value = 24664.413
$ws.Range("M136").Value = -2563.7586
$ws.Range("N136").Value = -29764.413

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 34773.4
$ws.Range("J70").Value = 34773.4
$ws.Range("L70").Value = 34773.4
$ws.Range("N70").Value = -35403.4
$ws.Range("H73").Value = 34773.4
$ws.Range("J73").Value = 34773.4
$ws.Range("L73").Value = 34773.4
$ws.Range("N73").Value = -36957.4
$ws.Range("H75").Value = 28963.3
$ws.Range("I75").Value = 29090
$ws.Range("J75").Value = 28949.223
$ws.Range("K75").Value = 29090
$ws.Range("L75").Value = 28949.223
$ws.Range("M75").Value = -28154
$ws.Range("N75").Value = -30821.223
$ws.Range("H78").Value = 28963.3
$ws.Range("I78").Value = 29090
$ws.Range("J78").Value = 28949.223
$ws.Range("K78").Value = 87270
$ws.Range("L78").Value = 86847.66900000001
$ws.Range("M78").Value = -82590
$ws.Range("N78").Value = -96207.66900000001
$ws.Range("H124").Value = 77249
$ws.Range("J124").Value = 77249
$ws.Range("K124").Value = 77249
$ws.Range("L124").Value = 77249
$ws.Range("N124").Value = -87069
